$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 24 ("More events") and Row 25 ("Light blending") gain a "Yes"
# marker in column B, matching the existing "Yes" cells (e.g. B7, B8,
# B21): same shared string value and the same highlighted cell style.
# Copy the format from an existing "Yes" cell, then stamp the value.

$ws.Range("B7").Copy() | Out-Null
$ws.Range("B24").PasteSpecial(-4122) | Out-Null
$ws.Range("B24").Value = "Yes"

$ws.Range("B7").Copy() | Out-Null
$ws.Range("B25").PasteSpecial(-4122) | Out-Null
$ws.Range("B25").Value = "Yes"

$excel.CutCopyMode = $false

# Move the active selection to B26, matching where editing left off.
$ws.Range("B26").Select() | Out-Null
